$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be forced to Text format
# first, otherwise Excel auto-converts the typed value into a numeric cell (the
# source workbook stores every Price/Volume cell as text).
$textCells = @("D5","D6","D7","D9","D10","D12","D13","D15","D17","D19","D20","D21","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '52.045.54'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '2.940.14'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').Value = '354.47'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').Value = '107.76'
$ws.Range('E6').Value = '  -4.52%  '
$ws.Range('D7').Value = '0.565'
$ws.Range('E7').Value = '  +1.66%  '
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('D9').Value = '0.622'
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').Value = '38.38'
$ws.Range('E10').Value = '  -3.58%  '
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').Value = '0.0866'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = '19.16'
$ws.Range('E13').Value = '  -3.17%  '
$ws.Range('D14').Value = '3.432.06'
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('D15').Value = '7.68'
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('D16').Value = '2.951.04'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '0.972'
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').Value = '52.091.02'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').Value = '3.46'
$ws.Range('E19').Value = '  +4.31%  '
$ws.Range('D20').Value = '7.50'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').Value = '13.62'
$ws.Range('E21').Value = '  -2.38%  '
$ws.Range('D22').Value = '0.0₃0976'
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').Value = '69.98'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D24').Value = '266.37'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').Value = '2.77'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').Value = '0.176'
$ws.Range('E26').Value = '  -3.78%  '
$ws.Range('D27').Value = '26.88'
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('D28').Value = '7.59'
$ws.Range('E28').Value = '  +13.87%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '0.106'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').Value = '10.31'
$ws.Range('E31').Value = '  -3.20%  '
$ws.Range('D32').Value = '36.69'
$ws.Range('E32').Value = '  -2.08%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').Value = '2.18'
$ws.Range('E33').Value = '  -3.45%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').Value = '6.02'
$ws.Range('E34').Value = '  -2.71%  '
$ws.Range('D35').Value = '52.21'
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('D36').Value = '0.0437'
$ws.Range('E36').Value = '  -2.80%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('D38').Value = '3.17'
$ws.Range('E38').Value = '  -4.19%  '
$ws.Range('D39').Value = '2.00'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('D40').Value = '17.89'
$ws.Range('E40').Value = '  -5.19%  '
$ws.Range('D41').Value = '2.69'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').Value = '0.118'
$ws.Range('E42').Value = '  +1.20%  '
$ws.Range('D43').Value = '22.93'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '118.67'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('D46').Value = '2.46'
$ws.Range('E46').Value = '  -4.38%  '
$ws.Range('D47').Value = '2.126.75'
$ws.Range('E47').Value = '  -2.03%  '
$ws.Range('D48').Value = '3.38'
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('D49').Value = '0.0345'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').Value = '0.240'
$ws.Range('E50').Value = '  -9.48%  '
$ws.Range('D51').Value = '1.35'
$ws.Range('E51').Value = '  +0.75%  '
